$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "ASMA_APT": update release date, and the two airports (Nice LFMN and
# Budapest LHBP) that now have a complete data set for 2021.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Release date (Meta data row) moves from 14 Feb 2022 to 14 Apr 2022.
$ws1.Range("B2").Value = 44665

# Airport name correction: Berlin/ Schoenefeld -> Berlin Brandenburg.
$ws1.Range("A7").Value = "Berlin Brandenburg (EDDB)"

# Row 30 = Marseille-Provence (LFML): new arrivals + additional ASMA time.
# E30 is a shared formula (=F30/D30) and recalculates automatically.
$ws1.Range("D30").Value = 26700
$ws1.Range("F30").Value = 14448

# Row 35 = Budapest/ Ferihegy (LHBP): new arrivals + additional ASMA time.
# E35 is a shared formula (=F35/D35) and recalculates automatically.
$ws1.Range("D35").Value = 25446
$ws1.Range("F35").Value = 17072
$ws1.Range("G35").Value = 13.83

# ---------------------------------------------------------------------------
# Sheet "Change Log": append the new change-log entry describing this update.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Copy the formatting from the existing date entry (A2) so the new date cell
# (A3) keeps the same "d mmm yyyy" style rather than getting a brand new xf.
$ws2.Range("A2").Copy()
$ws2.Range("A3").PasteSpecial(-4122)

$ws2.Range("A3").Value = 44665
$ws2.Range("B3").Value = "LFMN, LHBP"
$ws2.Range("C3").Value = 2021
$ws2.Range("D3").Value = "Airports updated with complete data set for 2021"
